$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look numeric (e.g. "214.67") need to be
# forced to remain plain text (matching the source inlineStr cells),
# otherwise Excel auto-converts them to numbers and changes the stored
# representation (e.g. dropping trailing zeros).
$textGuardCells = @(
    'D5', 'D8', 'D10', 'D14', 'D16', 'D19', 'D21', 'D25', 'D27', 'D30', 'D32', 'D33', 'D38', 'D41', 'D43', 'D46', 'D47'
)
foreach ($addr in $textGuardCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '27.098.09'
$ws.Range("E2").Value = '  -0.31%  '

# Row 3
$ws.Range("D3").Value = '1.624.02'
$ws.Range("E3").Value = '  -1.01%  '

# Row 4
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").Value = '214.67'
$ws.Range("E5").Value = '  -1.04%  '

# Row 6
$ws.Range("E6").Value = '  -1.11%  '

# Row 7
$ws.Range("E7").Value = '  -0.11%  '

# Row 8
$ws.Range("D8").Value = '0.0631'
$ws.Range("E8").Value = '  +0.89%  '

# Row 9
$ws.Range("E9").Value = '  -1.65%  '

# Row 10
$ws.Range("D10").Value = '20.09'
$ws.Range("E10").Value = '  +0.72%  '

# Row 12
$ws.Range("D12").Value = '1.852.01'
$ws.Range("E12").Value = '  -0.98%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.643.42'
$ws.Range("E13").Value = '  +0.11%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '4.15'
$ws.Range("E14").Value = '  +0.35%  '

# Row 15
$ws.Range("E15").Value = '  -0.09%  '

# Row 16
$ws.Range("D16").Value = '64.78'
$ws.Range("E16").Value = '  -3.28%  '

# Row 17
$ws.Range("D17").Value = '27.058.01'
$ws.Range("E17").Value = '  -0.49%  '

# Row 18
$ws.Range("E18").Value = '  +0.73%  '

# Row 19
$ws.Range("D19").Value = '214.25'
$ws.Range("E19").Value = '  -1.60%  '

# Row 20
$ws.Range("E20").Value = '  -0.15%  '

# Row 21
$ws.Range("D21").Value = '6.83'
$ws.Range("E21").Value = '  -1.30%  '

# Row 22
$ws.Range("E22").Value = '  -1.13%  '

# Row 23
$ws.Range("E23").Value = '  -6.74%  '

# Row 24
$ws.Range("E24").Value = '  -0.51%  '

# Row 25
$ws.Range("D25").Value = '148.07'
$ws.Range("E25").Value = '  +0.88%  '

# Row 26
$ws.Range("E26").Value = '  -0.14%  '

# Row 27
$ws.Range("D27").Value = '7.38'
$ws.Range("E27").Value = '  -0.77%  '

# Row 28
$ws.Range("E28").Value = '  -2.84%  '

# Row 29
$ws.Range("E29").Value = '  -0.49%  '

# Row 30
$ws.Range("D30").Value = '0.0512'
$ws.Range("E30").Value = '  +0.49%  '

# Row 31
$ws.Range("E31").Value = '  -1.03%  '

# Row 32
$ws.Range("D32").Value = '3.35'
$ws.Range("E32").Value = '  -0.89%  '

# Row 33
$ws.Range("D33").Value = '0.745'
$ws.Range("E33").Value = '  +36.35%  '

# Row 34
$ws.Range("E34").Value = '  -0.13%  '

# Row 35
$ws.Range("D35").Value = '1.358.21'
$ws.Range("E35").Value = '  +3.75%  '

# Row 36
$ws.Range("E36").Value = '  +0.51%  '

# Row 37
$ws.Range("E37").Value = '  -0.80%  '

# Row 38
$ws.Range("D38").Value = '0.0177'
$ws.Range("E38").Value = '  +1.05%  '

# Row 39
$ws.Range("E39").Value = '  -1.27%  '

# Row 41
$ws.Range("D41").Value = '0.804'
$ws.Range("E41").Value = '  -0.96%  '

# Row 43
$ws.Range("D43").Value = '65.06'
$ws.Range("E43").Value = '  +5.14%  '

# Row 44
$ws.Range("E44").Value = '  +1.10%  '

# Row 45
$ws.Range("D45").Value = '1.762.98'
$ws.Range("E45").Value = '  -1.00%  '

# Row 46
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '90.09'
$ws.Range("E46").Value = '  -1.62%  '

# Row 47
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = '0.880'
$ws.Range("E47").Value = '  +32.01%  '

# Row 48
$ws.Range("E48").Value = '  +2.79%  '

# Row 49
$ws.Range("E49").Value = '  -0.31%  '

# Row 50
$ws.Range("E50").Value = '  +5.57%  '

# Row 51
$ws.Range("E51").Value = '  +0.38%  '

# Restore the default cell style on the text-guarded cells so the
# workbook does not retain a stray explicit "@" number format.
foreach ($addr in $textGuardCells) {
    $ws.Range($addr).Style = "Normal"
}
